$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete rows 3 and 4 (no longer present in the updated layout) ---
$ws.Range("A4:J4").EntireRow.Delete()
$ws.Range("A3:J3").EntireRow.Delete()

# --- Insert the new "Tipo" column before the old column C (Tiene PT) ---
$ws.Range("C1:C2").EntireColumn.Insert()

# --- Insert two new columns after "Cantidad Neta" (now column G) for
#     "UND/ML/GR" and "Composicion de Lote", before the "Estado" column ---
$ws.Range("H1:H2").EntireColumn.Insert()
$ws.Range("H1:H2").EntireColumn.Insert()

# --- Row 1: headers ---
$ws.Cells.Item(1, 1).Value = "EAN/SKU/ID"
$ws.Cells.Item(1, 2).Value = "MARCA/TITULO"
$ws.Cells.Item(1, 3).Value = "Tipo"
$ws.Cells.Item(1, 4).Value = "Tiene PT"
$ws.Cells.Item(1, 5).Value = "Tiene ES"
$ws.Cells.Item(1, 6).Value = "Tiene IT"
$ws.Cells.Item(1, 7).Value = "Cantidad Neta"
$ws.Cells.Item(1, 8).Value = "UND/ML/GR"
$ws.Cells.Item(1, 9).Value = "Composición de Lote"
$ws.Cells.Item(1, 10).Value = "Estado"

# --- Row 2: data for the single remaining product ---
$ws.Cells.Item(2, 1).Value = "8809844997611"
$ws.Cells.Item(2, 2).Value = "Dr.Jart+ | Crema de Manos Hidratante | 100ml"
$ws.Cells.Item(2, 3).Value = "LOTE"
$ws.Cells.Item(2, 4).Value = "Tiene PT"
$ws.Cells.Item(2, 5).Value = "Tiene ES"
$ws.Cells.Item(2, 6).Value = "Tiene IT"
$ws.Cells.Item(2, 7).Value = "12"
$ws.Cells.Item(2, 8).Value = "UND"
$ws.Cells.Item(2, 9).Value = '"8809844997611","8809844997611","8809844997611","5245454545545","2323232332323"'
$ws.Cells.Item(2, 10).Value = "Solo Revisión"
